# Rename the two logo pictures that live in the document's headers/footers.
#
#   * Pearson logo (footers)  -> wp:docPr / pic:cNvPr  name: image2.png -> image1.png
#   * BTEC logo     (headers) -> wp:docPr / pic:cNvPr  name: image1.jpg -> image2.jpg
#
# InlineShape objects don't expose a settable "Name" (that's only on the
# floating Shape object in the Word OM), so each picture is briefly
# converted to a floating shape, renamed, then converted back to an
# inline shape (wdWrapInline) so the surrounding layout/XML shape is left
# exactly as it was.

$d = $word.ActiveDocument

function Rename-LogoPicture($shape, [string]$newName) {
    $floating = $shape.ConvertToShape()
    $floating.Name = $newName
    $floating.WrapFormat.Type = 7   # wdWrapInline - convert back to an inline picture
}

for ($s = 1; $s -le $d.Sections.Count; $s++) {
    $sec = $d.Sections.Item($s)

    # Headers: BTEC logo, image1.jpg -> image2.jpg
    for ($h = 1; $h -le 3; $h++) {
        $hdr = $sec.Headers.Item($h)
        if ($hdr.Exists) {
            for ($i = 1; $i -le $hdr.Range.InlineShapes.Count; $i++) {
                $shp = $hdr.Range.InlineShapes.Item($i)
                if ($shp.AlternativeText -eq "BTec_Logo-Orange") {
                    Rename-LogoPicture $shp "image2.jpg"
                }
            }
        }
    }

    # Footers: Pearson logo, image2.png -> image1.png
    for ($f = 1; $f -le 3; $f++) {
        $ftr = $sec.Footers.Item($f)
        if ($ftr.Exists) {
            for ($i = 1; $i -le $ftr.Range.InlineShapes.Count; $i++) {
                $shp = $ftr.Range.InlineShapes.Item($i)
                if ($shp.AlternativeText -like "*PearsonLogo.png") {
                    Rename-LogoPicture $shp "image1.png"
                }
            }
        }
    }
}
